$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "43.112.51"
$ws.Range("E2").Value = "  +4.90%  "
$ws.Range("D3").Value = "2.243.57"
$ws.Range("E3").Value = "  +3.09%  "
$ws.Range("E4").Value = "  +0.15%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "246.08"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +3.77%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "0.620"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  +1.03%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "76.13"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  +8.15%  "
$ws.Range("E8").Value = "  -0.08%  "
$ws.Range("E9").Value = "  +6.52%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "41.37"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  +3.23%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.0936"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  +0.80%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "55.43"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  +0.49%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "6.99"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  +3.48%  "
$ws.Range("E14").Value = "  +0.82%  "
$ws.Range("D15").Value = "2.557.40"
$ws.Range("E15").Value = "  +2.35%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "14.75"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  +5.97%  "
$ws.Range("B17").Value = "Polygon"
$ws.Range("C17").Value = "https://coinranking.com/coin/uW2tk-ILY0ii+polygon-matic"
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "0.818"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  +0.95%  "
$ws.Range("B18").Value = "WrappedEther"
$ws.Range("C18").Value = "https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth"
$ws.Range("D18").Value = "2.232.21"
$ws.Range("E18").Value = "  +2.72%  "
$ws.Range("D19").Value = "43.033.60"
$ws.Range("E19").Value = "  +4.98%  "
$ws.Range("E20").Value = "  +4.11%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "71.15"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  +1.04%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "6.01"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  +1.21%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "10.48"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  +3.46%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "2.22"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  +14.27%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "230.94"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  +2.42%  "
$ws.Range("E26").Value = "  -0.02%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "10.99"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  +1.52%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "3.37"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  -4.64%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "2.26"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  +2.06%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "2.22"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  +0.82%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "38.05"
$ws.Range("D31").Style = "Normal"
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "174.13"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  +4.20%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "20.37"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  +2.01%  "
$ws.Range("E34").Value = "  +3.76%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "5.42"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  +4.99%  "
$ws.Range("E36").Value = "  +1.49%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "0.111"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  +7.35%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "4.36"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  +5.84%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.0334"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  +17.07%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "13.19"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  +8.42%  "
$ws.Range("E41").Value = "  +3.29%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "5.59"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  +3.14%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.201"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  +5.40%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "60.41"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  +0.48%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "105.47"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  +7.45%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "8.63"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  +3.72%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "0.0995"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  +2.61%  "
$ws.Range("B48").Value = "ARBITRUM"
$ws.Range("C48").Value = "https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb"
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "1.11"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  +1.82%  "
$ws.Range("B49").Value = "WOONetwork"
$ws.Range("C49").Value = "https://coinranking.com/coin/k-J3YwacF+woonetwork-woo"
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "0.446"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  +20.88%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "2.31"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  +3.27%  "
$ws.Range("E51").Value = "  +1.45%  "
